$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Bitcoin wallet address (C1) and its font (Arial 10, black)
$ws.Range("C1").Value = "1HGu34ZGUn5QrjWf2an3xwYoGuFmYQMZoQ"
$ws.Range("C1").Font.Name = "Arial"
$ws.Range("C1").Font.Size = 10
$ws.Range("C1").Font.Color = 0

# Update bank card numbers (C2, D2, E2) to new server's values
$ws.Range("C2").Value = "4048415041393584"
$ws.Range("D2").Value = "4048415002254353"
$ws.Range("E2").Value = "4048415042013231"

# Update Litecoin wallet address (C4) and its font (Arial 10, black)
$ws.Range("C4").Value = "MBhv8ZvV1TGotH8BmiuEtHzYVCAZpLcr3U"
$ws.Range("C4").Font.Name = "Arial"
$ws.Range("C4").Font.Size = 10
$ws.Range("C4").Font.Color = 0

# Move the selection to C3
$ws.Range("C3").Select()
